{"js": "// Remove the \"Should implement screen reading facilities\" bullet point\n// (the paragraph immediately following \"Implement a fun and engaging way\n// to interact with the app through a non-standard gesture\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = \"Should implement screen reading facilities\";\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.trim() === target) {\n    para.delete();\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# Remove the \"Should implement screen reading facilities\" bullet point\n# (the paragraph immediately following \"Implement a fun and engaging way\n# to interact with the app through a non-standard gesture\").\n$d = $word.ActiveDocument\n$target = \"Should implement screen reading facilities\"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq $target) {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
